$wb = $excel.ActiveWorkbook

# --- Status text update: "Ready for handoff" -> "In Translation" -------
# Overview sheet tracks per-locale status in columns E (zh-cn) and F (de-de).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# Each locale detail sheet repeats the same status in its "Status" column (C).
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- Narrow the "Status" columns to fit the shorter text -----------------
# (Excel's ColumnWidth setter snaps to whole-pixel character steps, so this
# lands on the closest reachable width to the target ~13.41 chars.)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
